# Apply updates to Neg_Change and Pos_Change sheets per target diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New data for the "Neg_Change" worksheet (16 data rows, down from 18).
# ---------------------------------------------------------------------------
$negData = @(
    ,@("HEROMOTOCO", 4603, 4624.7, 4539, 4565, 455558, 972461, -0.5315411106460824, "HEROMOTOCO")
    ,@("GODREJCP", 1187.1, 1203.1, 1159.6, 1200, 940924, 2018271, -0.5337969975290732, "GODREJCP")
    ,@("PIDILITIND", 3099, 3101, 3065, 3095.9, 324415, 637714, -0.4912844943030888, "PIDILITIND")
    ,@("ICICIGI", 1894.9, 1918, 1882, 1900, 187046, 426758, -0.5617047600747965, "ICICIGI")
    ,@("OBEROIRLTY", 1567.4, 1612.5, 1564.3, 1611, 312265, 616455, -0.4934504546155032, "OBEROIRLTY")
    ,@("INDUSTOWER", 335, 338, 330.05, 337, 6020830, 13407199, -0.5509255885588034, "INDUSTOWER")
    ,@("KPITTECH", 1190.8, 1208.8, 1185.8, 1198.3, 628290, 1236355, -0.4918207149241116, "KPITTECH")
    ,@("BSE", 2400, 2421, 2370, 2399, 3272568, 6685163, -0.5104729682731745, "BSE")
    ,@("NATIONALUM", 188.4, 189.65, 185.73, 186.69, 4836690, 9906271, -0.511754725870108, "NATIONALUM")
    ,@("JSL", 677, 693.2, 674.2, 677.45, 771011, 1863449, -0.5862451829913241, "JSL")
    ,@("KFINTECH", 1073.1, 1097, 1068.4, 1093.8, 883832, 1959646, -0.5489838470825853, "KFINTECH")
    ,@("NCC", 221.95, 225.4, 220.33, 223.61, 2065238, 4371016, -0.5275153419708369, "NCC")
    ,@("PNBHOUSING", 774.65, 780, 766.25, 774, 1790243, 3559039, -0.4969869675493862, "PNBHOUSING")
    ,@("PPLPHARMA", 190.07, 191.11, 188, 189.89, 1608219, 3557025, -0.5478752609273199, "PPLPHARMA")
    ,@("CDSL", 1572, 1584, 1546.1, 1565, 2750081, 6086096, -0.5481370980674639, "CDSL")
    ,@("CROMPTON", 318.1, 320.3, 314, 317.5, 1767812, 3901406, -0.5468782280029302, "CROMPTON")
)

# ---------------------------------------------------------------------------
# New data for the "Pos_Change" worksheet (17 data rows, up from 15).
# ---------------------------------------------------------------------------
$posData = @(
    ,@("ETERNAL", 300.9, 310.4, 299.5, 309, 44948318, 28458176, 0.5794518243193099, "ETERNAL")
    ,@("APOLLOHOSP", 7101, 7320, 7020.5, 7248, 498598, 327102, 0.5242890596816895, "APOLLOHOSP")
    ,@("POWERGRID", 284.85, 286.4, 283.3, 285.4, 8755499, 6172968, 0.4183613133909004, "POWERGRID")
    ,@("ICICIBANK", 1427.9, 1438.4, 1420.4, 1437.3, 8492484, 5415511, 0.5681777767601247, "ICICIBANK")
    ,@("BAJFINANCE", 877.15, 888.3, 872.75, 876.45, 5367926, 3736916, 0.4364588339689733, "BAJFINANCE")
    ,@("PFC", 405.9, 419.9, 405.45, 419, 9713008, 6103086, 0.5914912554075102, "PFC")
    ,@("JSWENERGY", 521.5, 537, 518.3, 534.65, 2511889, 1688157, 0.4879475072519914, "JSWENERGY")
    ,@("TVSMOTOR", 2985, 2999, 2957.3, 2981.9, 569270, 387198, 0.4702297015996983, "TVSMOTOR")
    ,@("HAL", 4460, 4545, 4439, 4444, 1022011, 659530, 0.5496050217579185, "HAL")
    ,@("INDIANB", 650.85, 677, 650, 676.55, 2887066, 1934270, 0.4925868673970025, "INDIANB")
    ,@("NHPC", 82.45, 83.90000000000001, 81.77, 83.54000000000001, 10068007, 6613741, 0.5222862522133842, "NHPC")
    ,@("MAXHEALTH", 1245.3, 1270.6, 1243.3, 1261, 807096, 548971, 0.4701978793050999, "MAXHEALTH")
    ,@("SAIL", 121.07, 121.93, 119.9, 121.51, 9091136, 6142679, 0.47999529195649, "SAIL")
    ,@("IGL", 202.17, 203.59, 200.28, 202.5, 1404201, 995233, 0.4109268884773716, "IGL")
    ,@("HINDZINC", 420.7, 422.9, 414.55, 416.85, 1862022, 1210138, 0.5386856705598866, "HINDZINC")
    ,@("TATAELXSI", 5759.5, 5769.5, 5665, 5670, 101170, 66575, 0.5196395043184379, "TATAELXSI")
    ,@("NBCC", 105.69, 109.95, 104.61, 109.25, 12682149, 8350031, 0.5188146008080688, "NBCC")
)

$wsNeg = $wb.Worksheets.Item("Neg_Change")
$wsPos = $wb.Worksheets.Item("Pos_Change")

# --- Neg_Change: overwrite rows 2..17 with new values ----------------------
for ($i = 0; $i -lt $negData.Length; $i++) {
    $r = $i + 2
    $vals = $negData[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $wsNeg.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
# Remove the now-unused 18th row (sheet shrinks from 18 to 17 rows).
$wsNeg.Range("A18:I18").ClearContents()

# --- Pos_Change: overwrite rows 2..18 with new values -----------------------
for ($i = 0; $i -lt $posData.Length; $i++) {
    $r = $i + 2
    $vals = $posData[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $wsPos.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

Write-Output "Updated Neg_Change ($($negData.Length) rows) and Pos_Change ($($posData.Length) rows)"
